# Applies the update described by the diff:
#  - Sheet1 row 3: cells C3,F3,H3,I3,J3,K3,N3 change their text from
#    "departing on 11/09/2019" (shared string 90) to a brand-new shared
#    string "departing on 29/02/2020", and pick up the date-style format
#    (the same cell style already used by C3, numFmt "m/d/yyyy").
#  - Sheet1's active selection moves from C18 to P8.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$targets = @("C3", "F3", "H3", "I3", "J3", "K3", "N3")
$newText = "departing on 29/02/2020"

# 1) Reset the targeted cells to the plain/general style first (donor: A1,
#    which already carries the default "General" style). Doing this before
#    writing the new value stops Excel from re-deriving a brand new ad-hoc
#    number format out of the date-looking text while the old date format
#    is still applied (which would otherwise pollute styles.xml).
$ws1.Range("A1").Copy() | Out-Null
foreach ($addr in $targets) {
    $ws1.Range($addr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

# 2) Write the new text into each target cell.
foreach ($addr in $targets) {
    $ws1.Range($addr).Value = $newText
}

# 3) Re-apply the original date-style formatting (same style already used
#    by these cells), sourced from Sheet2!B3 which keeps that exact style
#    and is untouched by this change.
$ws2.Range("B3").Copy() | Out-Null
foreach ($addr in $targets) {
    $ws1.Range($addr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

# 4) Move the active selection on Sheet1 to P8.
$ws1.Range("P8").Select() | Out-Null
